$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = '27/12/2025 10:03'
$ws.Cells.Item(46, 3).Value = 163
$ws.Cells.Item(46, 4).Value = 'Português'
$ws.Cells.Item(46, 5).Value = 'Vírgula'
$ws.Cells.Item(46, 6).Value = 1

# Row 47
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = '27/12/2025 10:04'
$ws.Cells.Item(47, 3).Value = 70
$ws.Cells.Item(47, 4).Value = 'Português'
$ws.Cells.Item(47, 5).Value = 'Conjunção'
$ws.Cells.Item(47, 6).Value = 1

# Row 48
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = '27/12/2025 10:05'
$ws.Cells.Item(48, 3).Value = 108
$ws.Cells.Item(48, 4).Value = 'Português'
$ws.Cells.Item(48, 5).Value = 'Emprego de Tempos e Modos'
$ws.Cells.Item(48, 6).Value = 1

# Row 49
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = '27/12/2025 10:07'
$ws.Cells.Item(49, 3).Value = 157
$ws.Cells.Item(49, 4).Value = 'Português'
$ws.Cells.Item(49, 5).Value = 'Vírgula'
$ws.Cells.Item(49, 6).Value = 1

# Row 50
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = '27/12/2025 10:08'
$ws.Cells.Item(50, 3).Value = 76
$ws.Cells.Item(50, 4).Value = 'Português'
$ws.Cells.Item(50, 5).Value = 'Conjunção'
$ws.Cells.Item(50, 6).Value = 1

# Row 51
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = '27/12/2025 10:10'
$ws.Cells.Item(51, 3).Value = 112
$ws.Cells.Item(51, 4).Value = 'Português'
$ws.Cells.Item(51, 5).Value = 'Emprego de Tempos e Modos'
$ws.Cells.Item(51, 6).Value = 0

# Row 52
$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = '27/12/2025 10:13'
$ws.Cells.Item(52, 3).Value = 121
$ws.Cells.Item(52, 4).Value = 'Português'
$ws.Cells.Item(52, 5).Value = 'Verbos Traiçoeiros'
$ws.Cells.Item(52, 6).Value = 0

# Row 53
$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = '27/12/2025 10:15'
$ws.Cells.Item(53, 3).Value = 313
$ws.Cells.Item(53, 4).Value = 'Português'
$ws.Cells.Item(53, 5).Value = 'Narração'
$ws.Cells.Item(53, 6).Value = 1

# Row 54
$ws.Cells.Item(54, 1).Value = 53
$ws.Cells.Item(54, 2).Value = '27/12/2025 10:20'
$ws.Cells.Item(54, 3).NumberFormat = "@"
$ws.Cells.Item(54, 3).Value = '350'
$ws.Cells.Item(54, 3).Style = "Normal"
$ws.Cells.Item(54, 4).Value = 'Português'
$ws.Cells.Item(54, 5).Value = 'Compreensão E Interpretação'
$ws.Cells.Item(54, 6).Value = 1

# Row 55
$ws.Cells.Item(55, 1).Value = 54
$ws.Cells.Item(55, 2).Value = '27/12/2025 10:21'
$ws.Cells.Item(55, 3).NumberFormat = "@"
$ws.Cells.Item(55, 3).Value = '349'
$ws.Cells.Item(55, 3).Style = "Normal"
$ws.Cells.Item(55, 4).Value = 'Português'
$ws.Cells.Item(55, 5).Value = 'Compreensão E Interpretação'
$ws.Cells.Item(55, 6).Value = 1

# Row 56
$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 2).Value = '27/12/2025 10:21'
$ws.Cells.Item(56, 3).NumberFormat = "@"
$ws.Cells.Item(56, 3).Value = '348'
$ws.Cells.Item(56, 3).Style = "Normal"
$ws.Cells.Item(56, 4).Value = 'Português'
$ws.Cells.Item(56, 5).Value = 'Compreensão E Interpretação'
$ws.Cells.Item(56, 6).Value = 1

# Row 57
$ws.Cells.Item(57, 1).Value = 56
$ws.Cells.Item(57, 2).Value = '27/12/2025 10:23'
$ws.Cells.Item(57, 3).NumberFormat = "@"
$ws.Cells.Item(57, 3).Value = '181'
$ws.Cells.Item(57, 3).Style = "Normal"
$ws.Cells.Item(57, 4).Value = 'Português'
$ws.Cells.Item(57, 5).Value = 'Travessão'
$ws.Cells.Item(57, 6).Value = 0

# Row 58
$ws.Cells.Item(58, 1).Value = 57
$ws.Cells.Item(58, 2).Value = '27/12/2025 10:26'
$ws.Cells.Item(58, 3).NumberFormat = "@"
$ws.Cells.Item(58, 3).Value = '129'
$ws.Cells.Item(58, 3).Style = "Normal"
$ws.Cells.Item(58, 4).Value = 'Português'
$ws.Cells.Item(58, 5).Value = 'Funções Sintáticas'
$ws.Cells.Item(58, 6).Value = 1

# Row 59
$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(59, 2).Value = '27/12/2025 10:27'
$ws.Cells.Item(59, 3).Value = 168
$ws.Cells.Item(59, 4).Value = 'Português'
$ws.Cells.Item(59, 5).Value = 'Vírgula'
$ws.Cells.Item(59, 6).Value = 1

# Row 60
$ws.Cells.Item(60, 1).Value = 59
$ws.Cells.Item(60, 2).Value = '27/12/2025 10:30'
$ws.Cells.Item(60, 3).Value = 196
$ws.Cells.Item(60, 4).Value = 'Português'
$ws.Cells.Item(60, 5).Value = 'Concordância Verbal'
$ws.Cells.Item(60, 6).Value = 1

